# Add a new "Dried tangerine peel 陳皮" vocabulary entry to the word list,
# inserted right after the "As far as 據我所知" paragraph (and before the
# trailing empty/bookmark paragraph at the end of the document).

$d = $word.ActiveDocument

# Locate the paragraph that starts with "As far as" (i.e. "As far as 據我所知").
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "As far as*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'As far as ...' paragraph to insert after."
}

$targetPara = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new empty paragraph right after it.
$targetPara.Range.InsertParagraphAfter()

# Grab that freshly-created (still empty) paragraph...
$newPara = $d.Paragraphs.Item($targetIndex + 1)

# ...and populate it with the four runs that make up the new entry, using the
# WordOpenXML "pkg:package" envelope so each run keeps its own distinct
# formatting (the "tang" / "陳皮" runs carry an eastAsia font hint, and the
# "陳皮" run is additionally tagged zh-HK):
#   "Dried " | "tang" (rFonts hint=eastAsia) | "erine peel " | "陳皮" (rFonts hint=eastAsia, lang eastAsia=zh-HK)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/dummy.xml" pkg:contentType="application/xml">
    <pkg:xmlData>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:r><w:t xml:space="preserve">Dried </w:t></w:r>
        <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>tang</w:t></w:r>
        <w:r><w:t xml:space="preserve">erine peel </w:t></w:r>
        <w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-HK"/></w:rPr><w:t>陳皮</w:t></w:r>
      </w:p>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($xml) | Out-Null

Write-Output "Inserted 'Dried tangerine peel 陳皮' paragraph after paragraph $targetIndex."
